# Update "想去人数" (F column) counts on the 展览 sheet and the 全部类型 sheet.
# These two sheets list the same events (全部类型 also interleaves rows from
# 演出/本地生活), so the same set of F-value updates must be applied to both,
# using each sheet's own row numbering.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 116
$ws1.Range("F6").Value = 72
$ws1.Range("F7").Value = 2708
$ws1.Range("F9").Value = 1298
$ws1.Range("F12").Value = 10380
$ws1.Range("F15").Value = 276
$ws1.Range("F16").Value = 1020
$ws1.Range("F17").Value = 662
$ws1.Range("F18").Value = 11913
$ws1.Range("F19").Value = 12301

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 116
$ws4.Range("F6").Value = 72
$ws4.Range("F7").Value = 2708
$ws4.Range("F10").Value = 1298
$ws4.Range("F13").Value = 10380
$ws4.Range("F16").Value = 276
$ws4.Range("F17").Value = 1020
$ws4.Range("F18").Value = 662
$ws4.Range("F19").Value = 11913
$ws4.Range("F20").Value = 12301
